$wb = $excel.ActiveWorkbook
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
    Write-Host $s.Index
}
